$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (existing review) - email addresses changed
$ws.Range("C2").Value = "maorshmul94@gmail.com"
$ws.Range("D2").Value = "amramg25@gmail.com "

# Add new row 6 with a new review - copy formatting from the row above, then overwrite values
$ws.Range("A5:G5").Copy($ws.Range("A6:G6"))
$ws.Rows.Item(6).RowHeight = 13.8

$ws.Range("A6").Value = "com.hamxa.shaynachim"
$ws.Range("B6").Value = "bitcoin"
$ws.Range("C6").Value = "oamit1038@gmail.com"
$ws.Range("D6").Value = "maorshmul94@gmail.com"
$ws.Range("E6").Value = "27/5/2019 15:59"
$ws.Range("F6").Value = "remarkable info"
$ws.Range("G6").Value = "yes"

$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("D9").Select()
